$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 190.375
$ws.Range("I8").Value = 190.375
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 571.125
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -432.125
$ws.Range("N8").ClearContents()
$ws.Range("H17").Value = 2106.4443
$ws.Range("J17").Value = 2434.6667
$ws.Range("L17").Value = 7304.000100000001
$ws.Range("N17").Value = -7640.000100000001
$ws.Range("H19").Value = 162.25
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 183
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 183
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = -533
$ws.Range("H21").Value = 2187.25
$ws.Range("I21").Value = 1374.5
$ws.Range("J21").Value = 3000
$ws.Range("K21").Value = 1374.5
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = -906.5
$ws.Range("N21").Value = -3936
$ws.Range("H23").Value = 2187.25
$ws.Range("I23").Value = 1374.5
$ws.Range("J23").Value = 3000
$ws.Range("K23").Value = 1374.5
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = -1140.5
$ws.Range("N23").Value = -3468
$ws.Range("H28").Value = 114250.62
$ws.Range("I28").Value = 2545.5
$ws.Range("K28").Value = 2545.5
$ws.Range("M28").Value = -2060.5
$ws.Range("H70").Value = 1577.2106
$ws.Range("I70").Value = 1576.2142
$ws.Range("J70").Value = 1580
$ws.Range("K70").Value = 4728.642599999999
$ws.Range("L70").Value = 4740
$ws.Range("M70").Value = -4458.642599999999
$ws.Range("N70").Value = -5280
$ws.Range("H73").Value = 1577.2106
$ws.Range("I73").Value = 1576.2142
$ws.Range("J73").Value = 1580
$ws.Range("K73").Value = 4728.642599999999
$ws.Range("L73").Value = 4740
$ws.Range("M73").Value = -3792.642599999999
$ws.Range("N73").Value = -6612
$ws.Range("H82").Value = 2834.6667
$ws.Range("I82").Value = 2834.6667
$ws.Range("K82").Value = 8504.000100000001
$ws.Range("M82").Value = -8098.000100000001
$ws.Range("H85").Value = 2834.6667
$ws.Range("I85").Value = 2834.6667
$ws.Range("K85").Value = 8504.000100000001
$ws.Range("M85").Value = -7100.000100000001
$ws.Range("H86").Value = 208335680
$ws.Range("I86").Value = 83335864
$ws.Range("J86").Value = 333335500
$ws.Range("K86").Value = 83335864
$ws.Range("L86").Value = 333335500
$ws.Range("M86").Value = -83334741
$ws.Range("N86").Value = -333337746
$ws.Range("H89").Value = 208335680
$ws.Range("I89").Value = 83335864
$ws.Range("J89").Value = 333335500
$ws.Range("K89").Value = 416679320
$ws.Range("L89").Value = 1666677500
$ws.Range("M89").Value = -416673704
$ws.Range("N89").Value = -1666688732
$ws.Range("H92").Value = 254.6875
$ws.Range("I92").Value = 213.35715
$ws.Range("J92").Value = 544
$ws.Range("K92").Value = 213.35715
$ws.Range("L92").Value = 544
$ws.Range("M92").Value = 1034.64285
$ws.Range("N92").Value = -3040
$ws.Range("H98").Value = 1006.94116
$ws.Range("I98").Value = 882.375
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 882.375
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 615.625
$ws.Range("N98").Value = -5996
$ws.Range("H100").Value = 2909.2
$ws.Range("I100").Value = 2284.7144
$ws.Range("J100").Value = 4366.3335
$ws.Range("K100").Value = 2284.7144
$ws.Range("L100").Value = 4366.3335
$ws.Range("M100").Value = -1743.7144
$ws.Range("N100").Value = -5448.3335
$ws.Range("H106").Value = 279473.06
$ws.Range("I106").Value = 557526.8
$ws.Range("K106").Value = 557526.8
$ws.Range("M106").Value = -556895.8
$ws.Range("H112").Value = 1705.6364
$ws.Range("J112").Value = 1680.8125
$ws.Range("L112").Value = 5042.4375
$ws.Range("N112").Value = -7258.4375
$ws.Range("H116").Value = 2574033.8
$ws.Range("I116").Value = 8825.75
$ws.Range("K116").Value = 8825.75
$ws.Range("M116").Value = -5383.75
$ws.Range("H122").Value = 1006.94116
$ws.Range("I122").Value = 882.375
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2647.125
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -197.125
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 78048.07000000001
$ws.Range("I125").Value = 139842.62
$ws.Range("K125").Value = 1258583.58
$ws.Range("M125").Value = -1256123.58
$ws.Range("H129").Value = 1881
$ws.Range("J129").Value = 3317.4
$ws.Range("L129").Value = 9952.200000000001
$ws.Range("N129").Value = -19952.2
$ws.Range("H132").Value = 2684.5652
$ws.Range("I132").Value = 2345.3684
$ws.Range("J132").Value = 4295.75
$ws.Range("K132").Value = 7036.1052
$ws.Range("L132").Value = 12887.25
$ws.Range("M132").Value = -4506.1052
$ws.Range("N132").Value = -17947.25
$ws.Range("H137").Value = 277012.94
$ws.Range("I137").Value = 3198.2104
$ws.Range("J137").Value = 970676.9399999999
$ws.Range("K137").Value = 9594.6312
$ws.Range("L137").Value = 2912030.82
$ws.Range("M137").Value = -7044.6312
$ws.Range("N137").Value = -2917130.82
$ws.Range("H138").Value = 2383.738
$ws.Range("J138").Value = 2899.84
$ws.Range("L138").Value = 8699.52
$ws.Range("N138").Value = -18979.52

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8329.666999999999
$ws.Range("I8").Value = 5494.5
$ws.Range("K8").Value = 5494.5
$ws.Range("M8").Value = -5350.5
$ws.Range("H32").Value = 5584.269
$ws.Range("I32").Value = 3129.3374
$ws.Range("K32").Value = 3129.3374
$ws.Range("M32").Value = -2842.3374
$ws.Range("H45").Value = 6257244
$ws.Range("I45").Value = 9305.643
$ws.Range("K45").Value = 9305.643
$ws.Range("M45").Value = -8928.643
$ws.Range("H61").Value = 1943.738
$ws.Range("I61").Value = 1649.5
$ws.Range("J61").Value = 2885.3
$ws.Range("K61").Value = 1649.5
$ws.Range("L61").Value = 2885.3
$ws.Range("M61").Value = -1437.5
$ws.Range("N61").Value = -3309.3
$ws.Range("H74").Value = 2539.3157
$ws.Range("I74").Value = 1771.125
$ws.Range("K74").Value = 1771.125
$ws.Range("M74").Value = -897.125
$ws.Range("H77").Value = 2539.3157
$ws.Range("I77").Value = 1771.125
$ws.Range("K77").Value = 8855.625
$ws.Range("M77").Value = -4487.625
$ws.Range("H132").Value = 2225.6191
$ws.Range("I132").Value = 1675.9333
$ws.Range("J132").Value = 3599.8333
$ws.Range("K132").Value = 5027.7999
$ws.Range("L132").Value = 10799.4999
$ws.Range("M132").Value = -2497.7999
$ws.Range("N132").Value = -15859.4999
$ws.Range("H136").Value = 1943.738
$ws.Range("I136").Value = 1649.5
$ws.Range("J136").Value = 2885.3
$ws.Range("K136").Value = 4948.5
$ws.Range("L136").Value = 8655.900000000001
$ws.Range("M136").Value = -2398.5
$ws.Range("N136").Value = -13755.9

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 132816.78
$ws.Range("I20").Value = 192337.19
$ws.Range("K20").Value = 192337.19
$ws.Range("M20").Value = -192090.19
$ws.Range("H94").Value = 1467
$ws.Range("I94").Value = 1527.64
$ws.Range("J94").Value = 961.6667
$ws.Range("K94").Value = 1527.64
$ws.Range("L94").Value = 961.6667
$ws.Range("M94").Value = -1076.64
$ws.Range("N94").Value = -1863.6667
$ws.Range("H107").Value = 1719.174
$ws.Range("I107").Value = 1102.7858
$ws.Range("J107").Value = 2678
$ws.Range("K107").Value = 1102.7858
$ws.Range("L107").Value = 2678
$ws.Range("M107").Value = 817.2141999999999
$ws.Range("N107").Value = -6518
$ws.Range("H132").Value = 37729.973
$ws.Range("J132").Value = 37938.918
$ws.Range("L132").Value = 37938.918
$ws.Range("N132").Value = -48058.918
$ws.Range("H134").Value = 2205
$ws.Range("I134").Value = 1692.1428
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5076.428400000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2541.428400000001
$ws.Range("N134").Value = -17070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3346.65
$ws.Range("I31").Value = 2593.4167
$ws.Range("K31").Value = 2593.4167
$ws.Range("M31").Value = -2298.4167
$ws.Range("H34").Value = 3346.65
$ws.Range("I34").Value = 2593.4167
$ws.Range("K34").Value = 2593.4167
$ws.Range("M34").Value = -2391.4167
$ws.Range("H41").Value = 15812.375
$ws.Range("J41").Value = 37999.75
$ws.Range("L41").Value = 37999.75
$ws.Range("N41").Value = -38855.75
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -19434
$ws.Range("N47").ClearContents()
$ws.Range("H74").Value = 39999.332
$ws.Range("J74").Value = 39999.332
$ws.Range("L74").Value = 39999.332
$ws.Range("N74").Value = -41747.332
$ws.Range("H77").Value = 39999.332
$ws.Range("J77").Value = 39999.332
$ws.Range("L77").Value = 119997.996
$ws.Range("N77").Value = -128733.996
$ws.Range("H95").Value = 16666.666
$ws.Range("J95").Value = 16666.666
$ws.Range("L95").Value = 16666.666
$ws.Range("N95").Value = -22158.666
$ws.Range("H105").Value = 72349.69
$ws.Range("I105").Value = 125034
$ws.Range("J105").Value = 4612.7144
$ws.Range("K105").Value = 125034
$ws.Range("L105").Value = 4612.7144
$ws.Range("M105").Value = -123287
$ws.Range("N105").Value = -8106.7144
$ws.Range("H107").Value = 1755.2
$ws.Range("I107").Value = 1944
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1944
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -24
$ws.Range("N107").Value = -4840
$ws.Range("H132").Value = 1798.0244
$ws.Range("I132").Value = 1580.0294
$ws.Range("J132").Value = 2856.8572
$ws.Range("K132").Value = 4740.0882
$ws.Range("L132").Value = 8570.571599999999
$ws.Range("M132").Value = -2210.0882
$ws.Range("N132").Value = -13630.5716
$ws.Range("H141").Value = 132177.55
$ws.Range("J141").Value = 125465.7
$ws.Range("L141").Value = 125465.7
$ws.Range("N141").Value = -135825.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2595.5625
$ws.Range("J5").Value = 3535.6667
$ws.Range("L5").Value = 10607.0001
$ws.Range("N5").Value = -10831.0001
$ws.Range("H9").Value = 1083
$ws.Range("J9").Value = 1667
$ws.Range("L9").Value = 5001
$ws.Range("N9").Value = -5449
$ws.Range("H69").Value = 8753
$ws.Range("I69").Value = 10012
$ws.Range("J69").Value = 8333.333000000001
$ws.Range("K69").Value = 30036
$ws.Range("L69").Value = 24999.999
$ws.Range("M69").Value = -29225
$ws.Range("N69").Value = -26621.999
$ws.Range("H72").Value = 8753
$ws.Range("I72").Value = 10012
$ws.Range("J72").Value = 8333.333000000001
$ws.Range("K72").Value = 90108
$ws.Range("L72").Value = 74999.997
$ws.Range("M72").Value = -86052
$ws.Range("N72").Value = -83111.997
$ws.Range("H81").Value = 4932.6895
$ws.Range("J81").Value = 7374.9443
$ws.Range("L81").Value = 22124.8329
$ws.Range("N81").Value = -24370.8329
$ws.Range("H84").Value = 4932.6895
$ws.Range("J84").Value = 7374.9443
$ws.Range("L84").Value = 66374.4987
$ws.Range("N84").Value = -77606.4987
$ws.Range("H88").Value = 4823.8
$ws.Range("J88").Value = 4823.8
$ws.Range("L88").Value = 14471.4
$ws.Range("N88").Value = -15327.4
$ws.Range("H91").Value = 4823.8
$ws.Range("J91").Value = 4823.8
$ws.Range("L91").Value = 14471.4
$ws.Range("N91").Value = -17435.4
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242
$ws.Range("H130").Value = 29550
$ws.Range("I130").Value = 29550
$ws.Range("K130").Value = 88650
$ws.Range("M130").Value = -83630
$ws.Range("H133").Value = 8419.200000000001
$ws.Range("J133").Value = 9022
$ws.Range("L133").Value = 27066
$ws.Range("N133").Value = -37186
$ws.Range("H135").Value = 2595.5625
$ws.Range("J135").Value = 3535.6667
$ws.Range("L135").Value = 31821.0003
$ws.Range("N135").Value = -36891.0003
$ws.Range("H136").Value = 9000
$ws.Range("J136").Value = 9000
$ws.Range("L136").Value = 27000
$ws.Range("N136").Value = -37200
$ws.Range("H138").Value = 7185.533
$ws.Range("I138").Value = 7792.05
$ws.Range("J138").Value = 5972.5
$ws.Range("K138").Value = 23376.15
$ws.Range("L138").Value = 17917.5
$ws.Range("M138").Value = -18236.15
$ws.Range("N138").Value = -28197.5
$ws.Range("H139").Value = 2982.5789
$ws.Range("I139").Value = 1762.4286
$ws.Range("J139").Value = 6399
$ws.Range("K139").Value = 5287.2858
$ws.Range("L139").Value = 19197
$ws.Range("M139").Value = -147.2857999999997
$ws.Range("N139").Value = -29477
$ws.Range("H140").Value = 1592.4584
$ws.Range("I140").Value = 1522.5652
$ws.Range("J140").Value = 3200
$ws.Range("K140").Value = 4567.6956
$ws.Range("L140").Value = 9600
$ws.Range("M140").Value = 612.3044
$ws.Range("N140").Value = -19960

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1648.5714
$ws.Range("I2").Value = 2304
$ws.Range("K2").Value = 2304
$ws.Range("M2").Value = -2191
$ws.Range("H20").Value = 27490.428
$ws.Range("J20").Value = 29830.875
$ws.Range("L20").Value = 29830.875
$ws.Range("N20").Value = -30320.875
$ws.Range("H21").Value = 3335833.2
$ws.Range("I21").Value = 3335833.2
$ws.Range("K21").Value = 3335833.2
$ws.Range("M21").Value = -3335660.2
$ws.Range("H30").Value = 3335833.2
$ws.Range("I30").Value = 3335833.2
$ws.Range("K30").Value = 3335833.2
$ws.Range("M30").Value = -3335728.2
$ws.Range("H70").Value = 44151.27
$ws.Range("I70").Value = 47330.793
$ws.Range("K70").Value = 47330.793
$ws.Range("M70").Value = -47060.793
$ws.Range("H73").Value = 44151.27
$ws.Range("I73").Value = 47330.793
$ws.Range("K73").Value = 47330.793
$ws.Range("M73").Value = -46394.793
$ws.Range("H80").Value = 4442
$ws.Range("I80").Value = 4855.6665
$ws.Range("K80").Value = 4855.6665
$ws.Range("M80").Value = -3857.6665
$ws.Range("H83").Value = 4442
$ws.Range("I83").Value = 4855.6665
$ws.Range("K83").Value = 24278.3325
$ws.Range("M83").Value = -19286.3325
$ws.Range("H102").Value = 3583.5715
$ws.Range("I102").Value = 2866
$ws.Range("J102").Value = 4121.75
$ws.Range("K102").Value = 2866
$ws.Range("L102").Value = 4121.75
$ws.Range("M102").Value = -1244
$ws.Range("N102").Value = -7365.75
$ws.Range("H107").Value = 285.2
$ws.Range("I107").Value = 119
$ws.Range("J107").Value = 488.33334
$ws.Range("K107").Value = 119
$ws.Range("L107").Value = 488.33334
$ws.Range("M107").Value = 1801
$ws.Range("N107").Value = -4328.33334
$ws.Range("H126").Value = 2838.5625
$ws.Range("I126").Value = 2168.0833
$ws.Range("K126").Value = 6504.249899999999
$ws.Range("M126").Value = -4034.249899999999
$ws.Range("H132").Value = 5924.364
$ws.Range("I132").Value = 5355.7334
$ws.Range("J132").Value = 7142.857
$ws.Range("K132").Value = 16067.2002
$ws.Range("L132").Value = 21428.571
$ws.Range("M132").Value = -13537.2002
$ws.Range("N132").Value = -26488.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 25000
$ws.Range("I12").Value = 40000
$ws.Range("K12").Value = 40000
$ws.Range("M12").Value = -39830
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H46").Value = 6350.913
$ws.Range("I46").Value = 8465.4
$ws.Range("K46").Value = 8465.4
$ws.Range("M46").Value = -8277.4
$ws.Range("H61").Value = 5943.8887
$ws.Range("I61").Value = 5957.5
$ws.Range("K61").Value = 5957.5
$ws.Range("M61").Value = -5755.5
$ws.Range("H82").Value = 1481.1111
$ws.Range("I82").Value = 1525.7142
$ws.Range("J82").Value = 1325
$ws.Range("K82").Value = 1525.7142
$ws.Range("L82").Value = 1325
$ws.Range("M82").Value = -1164.7142
$ws.Range("N82").Value = -2047
$ws.Range("H85").Value = 1481.1111
$ws.Range("I85").Value = 1525.7142
$ws.Range("J85").Value = 1325
$ws.Range("K85").Value = 1525.7142
$ws.Range("L85").Value = 1325
$ws.Range("M85").Value = -277.7141999999999
$ws.Range("N85").Value = -3821
$ws.Range("H93").Value = 1561.8334
$ws.Range("I93").Value = 1561.8334
$ws.Range("K93").Value = 1561.8334
$ws.Range("M93").Value = -313.8334
$ws.Range("H113").Value = 5943.8887
$ws.Range("I113").Value = 5957.5
$ws.Range("K113").Value = 5957.5
$ws.Range("M113").Value = -3787.5
$ws.Range("H132").Value = 1984.1428
$ws.Range("I132").Value = 1644.4546
$ws.Range("J132").Value = 2357.8
$ws.Range("K132").Value = 4933.3638
$ws.Range("L132").Value = 7073.400000000001
$ws.Range("M132").Value = -2403.3638
$ws.Range("N132").Value = -12133.4
$ws.Range("H136").Value = 2622.077
$ws.Range("I136").Value = 2623.875
$ws.Range("J136").Value = 2619.2
$ws.Range("K136").Value = 7871.625
$ws.Range("L136").Value = 7857.599999999999
$ws.Range("M136").Value = -5321.625
$ws.Range("N136").Value = -12957.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 40000
$ws.Range("I3").Value = 40000
$ws.Range("K3").Value = 40000
$ws.Range("M3").Value = -39886
$ws.Range("H30").Value = 8200
$ws.Range("J30").Value = 8200
$ws.Range("L30").Value = 8200
$ws.Range("N30").Value = -8414
$ws.Range("H41").Value = 14656.625
$ws.Range("J41").Value = 15112.8
$ws.Range("L41").Value = 15112.8
$ws.Range("N41").Value = -15892.8
$ws.Range("H62").Value = 13469.4375
$ws.Range("J62").Value = 18097.75
$ws.Range("L62").Value = 18097.75
$ws.Range("N62").Value = -19345.75
$ws.Range("H65").Value = 13469.4375
$ws.Range("J65").Value = 18097.75
$ws.Range("L65").Value = 90488.75
$ws.Range("N65").Value = -96728.75
$ws.Range("H74").Value = 16749.834
$ws.Range("J74").Value = 15986
$ws.Range("L74").Value = 15986
$ws.Range("N74").Value = -17858
$ws.Range("H77").Value = 16749.834
$ws.Range("J77").Value = 15986
$ws.Range("L77").Value = 47958
$ws.Range("N77").Value = -57318
$ws.Range("H100").Value = 5495575
$ws.Range("I100").Value = 11905486
$ws.Range("J100").Value = 1365
$ws.Range("K100").Value = 23810972
$ws.Range("L100").Value = 2730
$ws.Range("M100").Value = -23810431
$ws.Range("N100").Value = -3812
$ws.Range("H107").Value = 4434.7837
$ws.Range("I107").Value = 5309.846
$ws.Range("J107").Value = 2366.4546
$ws.Range("K107").Value = 15929.538
$ws.Range("L107").Value = 7099.3638
$ws.Range("M107").Value = -14009.538
$ws.Range("N107").Value = -10939.3638
$ws.Range("H110").Value = 85000
$ws.Range("J110").Value = 85000
$ws.Range("L110").Value = 85000
$ws.Range("N110").Value = -93180
$ws.Range("H113").Value = 812
$ws.Range("I113").Value = 812
$ws.Range("K113").Value = 2436
$ws.Range("M113").Value = -266
$ws.Range("H122").Value = 3548.7407
$ws.Range("I122").Value = 3029.875
$ws.Range("K122").Value = 9089.625
$ws.Range("M122").Value = -6639.625
$ws.Range("H132").Value = 1369189
$ws.Range("I132").Value = 14530.954
$ws.Range("K132").Value = 43592.862
$ws.Range("M132").Value = -41062.862
$ws.Range("H136").Value = 1570.2222
$ws.Range("I136").Value = 1570.2222
$ws.Range("K136").Value = 4710.6666
$ws.Range("M136").Value = -2160.6666
